# The commit rewrites every paragraph of this OCR'd German translation
# (corrected/retranslated text, and several long lines split into
# shorter ones -- the paragraph count grows from 37 to 49), gives all
# but the first paragraph a small left indent, and flips the page to
# landscape. Because the text of essentially every paragraph changes,
# the cleanest reproduction is to re-author the whole body in one shot
# (via Range.InsertXML) rather than chain 40+ Find/Replace calls.

$d = $word.ActiveDocument

# Target paragraph texts, in order. Paragraph 1 keeps its (unstyled)
# look; every other paragraph gets <w:ind w:left="90"/>.
$paragraphText = @(
    "- BeeidigteUbersetzungaus dem Arabischen",
    "RepublikTunesien(Wappen dertunesischen Republik)",
    "MinisteriumfiirHochschulbildungundwissenschaftlicheForschung",
    "Universitit Manouba",
    "von",
    "FakultitflrLiteraturwissenschaft, Kinste undHumanwissenschaften",
    "Das Nationale Zeugnis (Bachelor)",
    "der Fundamentalen Lizenz",
    "7 Nach VoriagedesFilessesNr. des 1986vom01. |986 UberdieEmihtungderFabuliatder von Mancube,",
    "Voriagedes Nr 3deJtvesJahres2008 vor 25.SeptemberFebruat2008 wherdas Hochichulwesen,insbesondereLtersturwisserschaltdenenAnite",
    "Nat Gesecres 19 3.",
    "NashVecagedesEriassesNr 1932desJahces2008 vom .November 1992 oberdi Festlogung derStele,diefe neracichauagGerwisserschaftlichen nationalen Hochschulabychlesserstindi it,",
    "Nach Vertage desEriassesNr.3123 desJahres2008vom22.September208aberdieFestegungdesaligeracinen Rahracasde fordasStadiensystem und der BedsagungenfOrden Lrwer&gt; dex sexooalon",
    "HochschisbsehlusesfOrdieinsimdenderverschiedenenAusbiidungygrbicien,des Facher,StadiengtngenundFachrcltungen LMD-Sysiem",
    "Y ednachVorlagedecBeratangsprota.cllPrafungshommissionen im (1saenz.Masterokra",
    "Universitijahres2015-2016.",
    "wird Frau/ Fraulein: SalmaNjema (geboren 01. 01. 1993 inMonastir, Nationalausweisnummer:",
    "am 06935513)",
    "Das Nationale Zeugnis der Fundamentalen",
    "Lizenz (Bachelor) in: Fachbereich: Sprachen und Literaturen",
    "Hauptfach: Deutsche Sprache, Literatur und Landeskunde mit dem Pradikat: (Ausreichend) erteilt.",
    "Manouba,den02.07. 2016",
    "DerDekan: HabibKozdoghli(Unterschrift: Unleserlich) Dienstsiegel: (Fakultat frLiteraturwissenschaft,",
    "- Kansteund",
    "Humanwissenschaften-InderMitte: DerDekan).",
    "Trockenes Dienstsiegel: (Ministerium ftirHochschulbildungundwissenschafilicheForschungUniversitat von Manouba—InderMitte: Fakultiit",
    "fiirLiteraturwissenschaft,Kiinste und Humanwissenschaften Manouba)",
    "von",
    "Hinweis: Dasvorliegende Diplom wird nureinmalausgehiindigt.",
    "Aufder",
    "RUckseite:",
    "“Stempel des Ministeriums flirHochschulbildung und wissenschaftliche Forschung fllrdie Beglaubigung",
    "Durchsicht des Dokumentes: Beglaubigungsvermerk:",
    "erfolgte in derGeneraldirektion firHochschulbildung im MinisteriumfiirHochschulbildungundwissenschaftliche Forschung.",
    "wirdie EchtheitderUnterschrift Hiermitbestitigen",
    "des Herm: DerDekan ohneVerantwortung flrden Inhaltdesvorliegenden Dokuments, Beglaubigungsnummer: 3148, Ort",
    "Datum: Tunis, den 22. 0). 2025,",
    "u. Beglaubigungsgebiihr: 5 Dinar, Vizedirektorin der privaten Hochschulbildung: Latifa Ben Abderrahmen",
    "Unterschrift (Unleserlich), Siegel des Ministeriums fiir Hochschulbildung und wissenschafiliche Forschung (Republik Tunesien Ministerium flr",
    "HochschuibiidungundwissenschaftlicheForschung InderMitte: Wappender - tunesischen Republik)",
    "“Stempel des AuBenministeriums fir die Beglaubigung des Dokumentes: Beglaubiguagsvermerk: auswirtige",
    "Durchsicht erfolgte im Ministerium fir",
    "Angelegenheiten. Hiermit bestiitigen wir die Echtheit der Unterschrift der Frau: Latifa Bea Abderrahmen,",
    "i. A. des Ministers(Ur flir Hochschulhildung und",
    "wissenschaftliche Forschung,Ortund Datum: Tunis,den22. 01. 2025,Beglaubigungsgebiibr: 5 Dinar, i.A.des Ministers auswartige",
    "i. A. des Generaldirektors Angelegenheiten,",
    "fir konsularische Angelegenheiten: Hamida Labidi Unterschrift (Unleserlich), Siegel des Aufenministeriums (Republik",
    "Tunesien MinisteriumflrauswartigeAngelegenheiten ~",
    "- -InderMitte: Wappendertunesischen Republik) em"
)

$parts = @()
for ($i = 0; $i -lt $paragraphText.Count; $i++) {
    if ($i -eq 0) {
        $parts += "<w:p><w:r><w:t>" + $paragraphText[$i] + "</w:t></w:r></w:p>"
    } else {
        $parts += "<w:p><w:pPr><w:ind w:left=`"90`"/></w:pPr><w:r><w:t>" + $paragraphText[$i] + "</w:t></w:r></w:p>"
    }
}
$bodyXml = $parts -join ""

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole body in one shot (text + paragraph formatting).
[void]$d.Content.InsertXML($packageXml)

# Section: landscape Letter (swap width/height + orient="landscape").
$d.PageSetup.Orientation = 1

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
